$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 252.9574622683357
$ws.Range("G2").Value = 13.92717690194273
$ws.Range("H2").Value = 529.7361469707848
$ws.Range("I2").Value = 0.004060232748994343
$ws.Range("J2").Value = 0.00001236170850433409
$ws.Range("K2").Value = 0.01172832351857443
$ws.Range("L2").Value = 0.08717411495410761
$ws.Range("M2").Value = 0.00253898034227455
$ws.Range("N2").Value = 0.2021594616209604

# Row 3
$ws.Range("F3").Value = 0.0001704117928854918
$ws.Range("G3").Value = 0.00006625404822078588
$ws.Range("H3").Value = 0.0002844103565428984
$ws.Range("I3").Value = 0.0001580983872636318
$ws.Range("J3").Value = 0.0000618853884043025
$ws.Range("K3").Value = 0.0002624916828465526
$ws.Range("L3").Value = 0.0001765591310859227
$ws.Range("M3").Value = 0.00006854214567596484
$ws.Range("N3").Value = 0.0002947662782538867

# Row 4
$ws.Range("F4").Value = 252.9576326801287
$ws.Range("G4").Value = 13.92724315599094
$ws.Range("H4").Value = 529.7364313811413
$ws.Range("I4").Value = 0.004218331136257974
$ws.Range("J4").Value = 0.00007424709690863659
$ws.Range("K4").Value = 0.01199081520142098
$ws.Range("L4").Value = 0.08735067408519354
$ws.Range("M4").Value = 0.002607522487950515
$ws.Range("N4").Value = 0.2024542278992143
